# Weekly update: add a new day's price records (2 new rows) for
# Comercializadora del Agro de Limari - Limon.
# The new rows are inserted right before the existing block that starts
# at row 561 (date 44399), shifting that block (and everything after it)
# down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 561:562 (existing content shifts down).
$ws.Range("A561:A562").EntireRow.Insert()

# --- New row 561 ---
$ws.Range("A561").Value2 = 2
$ws.Range("B561").Value2 = "Comercializadora del Agro de Limarí"
$ws.Range("C561").Value2 = "Coquimbo"
$ws.Range("D561").Value2 = 44769
$ws.Range("E561").Value2 = 4
$ws.Range("F561").Value2 = "Fruta"
$ws.Range("G561").Value2 = 100102
$ws.Range("H561").Value2 = "Cítricos"
$ws.Range("I561").Value2 = 100102003
$ws.Range("J561").Value2 = "Limón"
$ws.Range("K561").Value2 = "Sin especificar"
$ws.Range("L561").Value2 = "1a amarillo"
$ws.Range("M561").Value2 = 750
$ws.Range("N561").Value2 = 2300
$ws.Range("O561").Value2 = 2500
$ws.Range("P561").Value2 = 2400
$ws.Range("Q561").Value2 = "`$/malla 16 kilos"
$ws.Range("R561").Value2 = "Provincia de Limarí"
$ws.Range("S561").Value2 = 150
$ws.Range("T561").Value2 = 16

# --- New row 562 ---
$ws.Range("A562").Value2 = 2
$ws.Range("B562").Value2 = "Comercializadora del Agro de Limarí"
$ws.Range("C562").Value2 = "Coquimbo"
$ws.Range("D562").Value2 = 44769
$ws.Range("E562").Value2 = 4
$ws.Range("F562").Value2 = "Fruta"
$ws.Range("G562").Value2 = 100102
$ws.Range("H562").Value2 = "Cítricos"
$ws.Range("I562").Value2 = 100102003
$ws.Range("J562").Value2 = "Limón"
$ws.Range("K562").Value2 = "Sin especificar"
$ws.Range("L562").Value2 = "2a amarillo"
$ws.Range("M562").Value2 = 600
$ws.Range("N562").Value2 = 1300
$ws.Range("O562").Value2 = 1500
$ws.Range("P562").Value2 = 1400
$ws.Range("Q562").Value2 = "`$/malla 16 kilos"
$ws.Range("R562").Value2 = "Provincia de Limarí"
$ws.Range("S562").Value2 = 88
$ws.Range("T562").Value2 = 16
